$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '60.819.42'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.18%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.366.72'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.52%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '568.80'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '137.63'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.68%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.74%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '7.65'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('E10').Value = '  -2.54%  '
$ws.Range('E11').Value = '  -4.79%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '3.944.97'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('E13').Value = '  +1.35%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.65'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.66%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.369.09'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('E16').Value = '  -2.31%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '60.977.32'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('E18').Value = '  -2.87%  '
$ws.Range('E19').Value = '  -3.99%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '8.88'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.65%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '381.05'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.94%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '75.71'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.84%  '
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('E24').Value = '  -0.06%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.0000109'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -6.95%  '
$ws.Range('E26').Value = '  +6.86%  '
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  -4.39%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.81'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  -7.19%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '22.88'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.42%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '167.65'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.86%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '6.81'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('E36').Value = '  -2.27%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.404.31'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('E38').Value = '  -3.65%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0755'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.87%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '25.17'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -9.62%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.770'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('E43').Value = '  -3.76%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.456.77'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.09'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('E46').Value = '  +0.07%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '6.60'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -3.81%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '22.05'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -6.85%  '
$ws.Range('E49').Value = '  -5.09%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.97'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -4.86%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.200'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.96%  '
